# Refresh the cryptos list (prices + 1h volume %, plus the TRON/WrappedEther
# rows swapping rank) to match the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ref, $text) {
    $cell = $ws.Range($ref)
    # Price cells hold digit-grouped text (e.g. "1.001", "26.302.01") that
    # Excel would otherwise silently reinterpret as a number; when the new
    # value still parses as a plain number, force text storage first.
    if ($text -match "^[+-]?[0-9]+(\.[0-9]+)?$") {
        $cell.NumberFormat = "@"
    }
    $cell.Value2 = $text
}

Set-CellText 'D2' '26.323.43'
Set-CellText 'E2' '  -3.05%  '
Set-CellText 'D3' '1.830.56'
Set-CellText 'E3' '  -2.74%  '
Set-CellText 'D4' '1.000'
Set-CellText 'E4' '  +0.04%  '
Set-CellText 'D5' '258.75'
Set-CellText 'E5' '  -8.16%  '
Set-CellText 'E6' '  +0.08%  '
Set-CellText 'D7' '0.5196'
Set-CellText 'E7' '  -1.92%  '
Set-CellText 'D8' '0.3225'
Set-CellText 'E8' '  -8.85%  '
Set-CellText 'D9' '0.06726'
Set-CellText 'E9' '  -4.58%  '
Set-CellText 'D10' '18.64'
Set-CellText 'E10' '  -8.70%  '
Set-CellText 'D11' '0.7643'
Set-CellText 'E11' '  -7.17%  '
Set-CellText 'B12' 'WrappedEther'
Set-CellText 'C12' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-CellText 'D12' '1.895.10'
Set-CellText 'E12' '  -0.78%  '
Set-CellText 'B13' 'TRON'
Set-CellText 'C13' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-CellText 'D13' '0.07676'
Set-CellText 'E13' '  -1.85%  '
Set-CellText 'D14' '88.60'
Set-CellText 'E14' '  -2.51%  '
Set-CellText 'D15' '5.015'
Set-CellText 'E15' '  -3.79%  '
Set-CellText 'D16' '0.9998'
Set-CellText 'E16' '  +0.05%  '
Set-CellText 'D17' '14.04'
Set-CellText 'E17' '  -3.94%  '
Set-CellText 'D18' '1.001'
Set-CellText 'E18' '  +0.12%  '
Set-CellText 'D19' '0.000007882'
Set-CellText 'E19' '  -3.72%  '
Set-CellText 'D20' '26.360.80'
Set-CellText 'E20' '  -3.08%  '
Set-CellText 'D21' '2.081.62'
Set-CellText 'E21' '  -1.84%  '
Set-CellText 'E22' '  -5.36%  '
Set-CellText 'D23' '9.397'
Set-CellText 'E23' '  -7.50%  '
Set-CellText 'D24' '5.887'
Set-CellText 'E24' '  -5.82%  '
Set-CellText 'D25' '2.283'
Set-CellText 'E25' '  -5.33%  '
Set-CellText 'D26' '145.35'
Set-CellText 'E26' '  -1.27%  '
Set-CellText 'D27' '1.640'
Set-CellText 'E27' '  -2.16%  '
Set-CellText 'D28' '16.91'
Set-CellText 'D29' '110.84'
Set-CellText 'E29' '  -3.16%  '
Set-CellText 'D30' '4.180'
Set-CellText 'E30' '  -5.85%  '
Set-CellText 'D31' '4.114'
Set-CellText 'E31' '  -6.48%  '
Set-CellText 'D32' '0.08715'
Set-CellText 'E32' '  -2.78%  '
Set-CellText 'D33' '0.04825'
Set-CellText 'E33' '  -2.56%  '
Set-CellText 'D34' '1.124'
Set-CellText 'E34' '  -5.15%  '
Set-CellText 'D35' '2.844'
Set-CellText 'E35' '  -2.11%  '
Set-CellText 'D36' '0.6815'
Set-CellText 'E36' '  -9.08%  '
Set-CellText 'D37' '3.084'
Set-CellText 'E37' '  -6.97%  '
Set-CellText 'D38' '0.01772'
Set-CellText 'D39' '2.205'
Set-CellText 'E39' '  -9.17%  '
Set-CellText 'D40' '0.4890'
Set-CellText 'E40' '  -8.36%  '
Set-CellText 'D41' '112.27'
Set-CellText 'E41' '  -4.23%  '
Set-CellText 'D42' '6.112'
Set-CellText 'E42' '  -3.64%  '
Set-CellText 'D43' '0.8854'
Set-CellText 'E43' '  -9.04%  '
Set-CellText 'E44' '  +0.17%  '
Set-CellText 'D45' '7.680'
Set-CellText 'E45' '  -6.92%  '
Set-CellText 'D46' '0.4190'
Set-CellText 'E46' '  -9.45%  '
Set-CellText 'D47' '0.1251'
Set-CellText 'E47' '  -8.84%  '
Set-CellText 'D48' '0.05874'
Set-CellText 'E48' '  -1.41%  '
Set-CellText 'D49' '9.015'
Set-CellText 'E49' '  -5.26%  '
Set-CellText 'D50' '35.26'
Set-CellText 'E50' '  -4.21%  '
Set-CellText 'D51' '59.33'
Set-CellText 'E51' '  -4.21%  '

Write-Output "Updated 99 cells"
